$d = $word.ActiveDocument

# This document has a single section whose header/footer pair is split
# across a "primary" (index 1) and a "first page" (index 2) story - that
# is why the same two logo pictures physically live in header1.xml /
# header2.xml and footer1.xml / footer2.xml. Walk both header and footer
# stories of every section and rename the inline picture(s) found there:
#   - the BTEC logo ("BTec_Logo-Orange")        image2.jpg -> image1.jpg
#   - the Pearson logo (PearsonLogo.png descr)  image1.png -> image2.png

for ($secIdx = 1; $secIdx -le $d.Sections.Count; $secIdx++) {
    $section = $d.Sections.Item($secIdx)

    for ($hfIdx = 1; $hfIdx -le 3; $hfIdx++) {
        $header = $section.Headers.Item($hfIdx)
        if ($header.Exists) {
            $shapes = $header.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shape = $shapes.Item($i)
                if ($shape.AlternativeText -eq "BTec_Logo-Orange") {
                    $shape.Name = "image1.jpg"
                } elseif ($shape.AlternativeText -like "*PearsonLogo.png") {
                    $shape.Name = "image2.png"
                }
            }
        }

        $footer = $section.Footers.Item($hfIdx)
        if ($footer.Exists) {
            $shapes = $footer.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shape = $shapes.Item($i)
                if ($shape.AlternativeText -eq "BTec_Logo-Orange") {
                    $shape.Name = "image1.jpg"
                } elseif ($shape.AlternativeText -like "*PearsonLogo.png") {
                    $shape.Name = "image2.png"
                }
            }
        }
    }
}
